# Jeceaba Planos workbook update
# - Updates the "Vibração" sheet's data table (new Hh/periodicidade figures,
#   one fewer row) and strips the heavy borders/number formats that table
#   used to carry.
# - Inserts a new blank worksheet "Planilha1" between "Vibração" and "Atual".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rewrite the "Vibração" table with the new figures
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vibração")

# Drop the last existing data row (old sheet had 18 data rows / rows 2-19,
# the new one only has 17 / rows 2-18).
$ws.Rows.Item(19).Delete()

# row, Atividade, singleRowGroup(A-style), Periodicidade, firstRowOfGroup(B-style), Hh
$rows = @(
    @(2,  "CR",  $false, 2, $true,  6),
    @(3,  "CR",  $false, 1, $false, 11.7),
    @(4,  "CR",  $false, 6, $false, 0.9),
    @(5,  "CR",  $false, 3, $false, 24.9),
    @(6,  "CV",  $false, 2, $true,  22.36),
    @(7,  "CV",  $false, 1, $false, 19.35),
    @(8,  "CV",  $false, 4, $false, 3.44),
    @(9,  "CV",  $false, 6, $false, 3.87),
    @(10, "CV",  $false, 3, $false, 18.92),
    @(11, "ENE", $true,  1, $true,  8.82),
    @(12, "LA",  $false, 1, $false, 28.6),
    @(13, "LA",  $false, 3, $false, 37.4),
    @(14, "NDT", $true,  1, $true,  11.22),
    @(15, "OFC", $true,  2, $true,  8.72),
    @(16, "PBL", $true,  2, $true,  3.12),
    @(17, "RK",  $false, 2, $true,  21.6),
    @(18, "RK",  $false, 1, $false, 52.65)
)

foreach ($row in $rows) {
    $r   = $row[0]
    $act = $row[1]
    $per = $row[3]
    $hh  = $row[5]

    $ws.Cells.Item($r, 1).Value = $act
    $ws.Cells.Item($r, 2).Value = $per
    $ws.Cells.Item($r, 3).Value = $hh
    $ws.Cells.Item($r, 4).Value = 1
}

# Reset formatting across the whole table body, then reapply only the
# font/fill/alignment that the refreshed table keeps (no borders, no
# dedicated number formats any more).
$body = $ws.Range("A2:D18")
$body.ClearFormats()

# Column A ("Atividade") - bold dark-grey Arial on white, top aligned.
$colA = $ws.Range("A2:A18")
$colA.Font.Name = "Arial"
$colA.Font.Bold = $true
$colA.Font.Size = 11
$colA.Font.Color = 2171169
$colA.Interior.Color = 16777215
$colA.VerticalAlignment = -4160

# Column B ("Periodicidade em meses") - bold dark-grey Arial on white,
# centered, wrapped.
$colB = $ws.Range("B2:B18")
$colB.Font.Name = "Arial"
$colB.Font.Bold = $true
$colB.Font.Size = 11
$colB.Font.Color = 2171169
$colB.Interior.Color = 16777215
$colB.HorizontalAlignment = -4108
$colB.WrapText = $true

# First row of each activity group is top-aligned instead of center-aligned.
foreach ($row in $rows) {
    if ($row[4]) {
        $ws.Cells.Item($row[0], 2).VerticalAlignment = -4160
    } else {
        $ws.Cells.Item($row[0], 2).VerticalAlignment = -4108
    }
}

# Column C ("Hh") - regular dark-grey Arial on white, right aligned, wrapped.
$colC = $ws.Range("C2:C18")
$colC.Font.Name = "Arial"
$colC.Font.Bold = $false
$colC.Font.Size = 11
$colC.Font.Color = 2171169
$colC.Interior.Color = 16777215
$colC.HorizontalAlignment = -4152
$colC.VerticalAlignment = -4108
$colC.WrapText = $true

# Column D ("Inicio") keeps the plain default look (already cleared above).

$ws.Range("S10").Select()

# ---------------------------------------------------------------------------
# 2) Insert the new blank "Planilha1" sheet right after "Vibração"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("Vibração")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "Planilha1"
$newSheet.Range("B11").Select()
